$wb = $excel.ActiveWorkbook

# --- Insert a new "Player Info" worksheet before the existing "ODI Batting" sheet ---
$battingSheetForInsert = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($battingSheetForInsert)
$newSheet.Name = "Player Info"

# Re-fetch stable references by name (positional refs shift after Add/rename)
$info = $wb.Worksheets.Item("Player Info")
$batting = $wb.Worksheets.Item("ODI Batting")

# Header row values
$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

# Copy the header formatting (bold, centered, bordered) from the ODI Batting header row
$batting.Range("A1:D1").Copy()
$info.Range("A1:D1").PasteSpecial(-4122)

# Data row - keep values as text (matches source workbook convention of storing all values as text)
$info.Range("A2:D2").NumberFormat = "@"
$info.Range("A2").Value = "7118"
$info.Range("B2").Value = "Muthuthanthirige Nuwanidu Keshawa Fernando"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Right Arm Off Break"

# --- Update the "ODI Batting" sheet ---
# Rename MATCH_CARD_LINK header to MATCH_CODE
$batting.Range("D1").Value = "MATCH_CODE"

# Replace the full scorecard URLs with just the numeric match code, keep as text
$batting.Range("D2:D4").NumberFormat = "@"
$batting.Range("D2").Value = "4689"
$batting.Range("D3").Value = "4691"
$batting.Range("D4").Value = "4735"

# Append a new 5th row of match data
$batting.Range("A5:J5").NumberFormat = "@"
$batting.Range("A5").Value = "4"
$batting.Range("B5").Value = "4"
$batting.Range("C5").Value = "31/03/2023"
$batting.Range("D5").Value = "4745"
$batting.Range("E5").Value = "1st"
$batting.Range("F5").Value = "New Zealand"
$batting.Range("G5").Value = "Seddon Park"
$batting.Range("H5").Value = "c H M Nicholls b M J Henry"
$batting.Range("I5").Value = "2"
$batting.Range("J5").Value = "9"

Write-Host "Edit complete"
